$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 86 (shifts rows 86:93 down to 87:94)
$ws.Rows.Item(86).Insert()

# Populate the newly inserted row 86 with the new weekly record
$ws.Range("A86").Value = 6
$ws.Range("B86").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C86").Value = "Metropolitana"
$ws.Range("D86").Value = 44776
$ws.Range("E86").Value = 13
$ws.Range("F86").Value = 100114007
$ws.Range("G86").Value = "Jengibre"
$ws.Range("H86").Value = "Sin especificar"
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 560
$ws.Range("K86").Value = 10000
$ws.Range("L86").Value = 12000
$ws.Range("M86").Value = 10857
$ws.Range("N86").Value = "$/caja 13 kilos"
$ws.Range("O86").Value = "Perú"
$ws.Range("P86").Value = 835
$ws.Range("Q86").Value = 13
$ws.Range("R86").Value = "Hortaliza"
